$d = $word.ActiveDocument

# The document ends with a paragraph that holds only the _GoBack bookmark
# (right after "I've attached a couple of Manhattan plots ..."). We need to
# insert, ahead of that bookmark paragraph:
#   1. a new blank paragraph
#   2. a new paragraph asking whether Gemma and bigRR agree on gene hits

$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertParagraphBefore()

$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertParagraphBefore()

$n = $d.Paragraphs.Count

# Paragraph n-2 is the first (blank) inserted paragraph -> make it truly empty.
$blankPara = $d.Paragraphs.Item($n - 2)
[void]$blankPara.Range.Duplicate.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# Paragraph n-1 is the second inserted paragraph -> fill it with the question,
# split across three runs (the middle one wrapped in proofErr spell-check
# markers around "bigRR"), all sharing the same run formatting.
$textPara = $d.Paragraphs.Item($n - 1)
$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="1F497D"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>'
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' `
     + '<w:r>' + $rPr + '<w:t xml:space="preserve">Do Gemma and </w:t></w:r>' `
     + '<w:proofErr w:type="spellStart"/>' `
     + '<w:r>' + $rPr + '<w:t>bigRR</w:t></w:r>' `
     + '<w:proofErr w:type="spellEnd"/>' `
     + '<w:r>' + $rPr + '<w:t xml:space="preserve"> find the same genes for the same traits?</w:t></w:r>' `
     + '</w:p>'
[void]$textPara.Range.Duplicate.InsertXML($xml)
